$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# --- Extend the header (row 1) and data rows (2-3) from columns B:G out to N,
#     and turn row 1 from a (wrong) duplicate of row 2's data into real column
#     headers: species, debtor, owner, total, register_date, register_reason,
#     property_category, category, date, legislator_name, legislator_id,
#     source_file, index ---

# Copy the existing header formatting (bold, bordered, centered) onto the new
# header cells H1:N1 before writing their labels.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 1 - column headers
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# The "date" column (J) holds a literal yyyy-mm-dd text value, not a real
# date - format it as text first so it isn't auto-converted to a date serial.
$ws.Range("J2:J3").NumberFormat = "@"

# Row 2 - first debt record (土地抵押權 / 華泰商業銀行...)
$ws.Range("B2").Value = "土地抵押權"
$ws.Range("C2").Value = "薛凌"
$ws.Range("D2").Value = "華泰商業銀行臺北市中山區長安東路"
$ws.Range("E2").Value = 95000000
$ws.Range("F2").Value = "99年11月05日"
$ws.Range("G2").Value = "抵押借款"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-03-30"
$ws.Range("K2").Value = "薛凌"
$ws.Range("L2").Value = 1384
$ws.Range("M2").Value = "tmpe95e1"
$ws.Range("N2").Value = 240

# Row 3 - second debt record (房貸 / 合作金庫商業銀行...)
$ws.Range("B3").Value = "房貸"
$ws.Range("C3").Value = "陳勝宏"
$ws.Range("D3").Value = "合作金庫商業銀行臺北市大安區和平東路"
$ws.Range("E3").Value = 31800000
$ws.Range("F3").Value = "99年06月30日"
$ws.Range("G3").Value = "房貸"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-03-30"
$ws.Range("K3").Value = "薛凌"
$ws.Range("L3").Value = 1384
$ws.Range("M3").Value = "tmpe95e1"
$ws.Range("N3").Value = 241
